$wb = $excel.ActiveWorkbook

# --- keyword_set sheet: rename "HR&L" keyword and add three new species keywords ---
$kwSheet = $wb.Worksheets.Item("keyword_set")

$kwSheet.Columns.Item(1).ColumnWidth = 24.5

$kwSheet.Range("A12").Value = "tule perch"
$kwSheet.Range("A13").Value = "speckled dace"
$kwSheet.Range("A14").Value = "sacramento pikeminnow"
$kwSheet.Range("A11").Value = "Healthy Rivers and Landscapes"

[void]$kwSheet.Range("F22").Select()

# --- funding sheet: add funder row and make it the active sheet ---
$fundingSheet = $wb.Worksheets.Item("funding")
$fundingSheet.Range("A2").Value = "California Department of Water Resources"

[void]$fundingSheet.Activate()
[void]$fundingSheet.Range("M23").Select()
